# Update column G (K = strike count) values for rows 2-81.
# The new values were recalculated when the save_data regeneration switched
# from using "Strike#" to "K", so here we just write the recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 3
    8  = 0
    9  = 2
    10 = 2
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 4
    17 = 1
    18 = 0
    19 = 4
    20 = 1
    21 = 0
    22 = 3
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 0
    29 = 2
    30 = 2
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 2
    39 = 2
    40 = 0
    41 = 1
    42 = 1
    43 = 2
    44 = 0
    45 = 1
    46 = 0
    47 = 0
    48 = 1
    49 = 2
    50 = 1
    51 = 1
    52 = 0
    53 = 2
    54 = 1
    55 = 2
    56 = 0
    57 = 1
    58 = 0
    59 = 0
    60 = 0
    61 = 4
    62 = 1
    63 = 0
    64 = 2
    65 = 1
    66 = 1
    67 = 2
    68 = 1
    69 = 0
    70 = 1
    71 = 1
    72 = 2
    73 = 0
    74 = 1
    75 = 0
    76 = 0
    77 = 1
    78 = 1
    79 = 1
    80 = 1
    81 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
